$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to match the latest scrape.
# Column D ("Price") values are forced to Text format before assignment so that
# numeric-looking strings (e.g. "231.79") are not auto-converted into Excel numbers,
# matching the original inline-string cell type. ClearFormats() afterwards removes
# the quote-prefix style marker so no stray style index is introduced.
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "42.528.72"
$c.ClearFormats()
$ws.Cells.Item(2, 5).Value = "  +2.56%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "2.230.25"
$c.ClearFormats()
$ws.Cells.Item(3, 5).Value = "  +0.48%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "231.79"
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +0.24%  "
$ws.Cells.Item(6, 5).Value = "  -0.64%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "60.92"
$c.ClearFormats()
$ws.Cells.Item(7, 5).Value = "  -0.20%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.406"
$c.ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +0.88%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.0908"
$c.ClearFormats()
$ws.Cells.Item(10, 5).Value = "  +2.24%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.ClearFormats()
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "2.560.97"
$c.ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +0.46%  "
$ws.Cells.Item(13, 5).Value = "  -0.88%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "22.29"
$c.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  +2.40%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "5.66"
$c.ClearFormats()
$ws.Cells.Item(15, 5).Value = "  +2.18%  "
$ws.Cells.Item(16, 5).Value = "  +0.13%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "2.234.26"
$c.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  +0.81%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "42.356.26"
$c.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  +2.52%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0943"
$c.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +5.24%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "6.18"
$c.ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +2.17%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "72.32"
$c.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  -0.85%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "244.96"
$c.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -2.22%  "
$ws.Cells.Item(23, 5).Value = "  -0.12%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.53"
$c.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +6.03%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.39"
$c.ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +5.48%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "9.73"
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  +1.84%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "169.52"
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +0.79%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "0.143"
$c.ClearFormats()
$ws.Cells.Item(28, 5).Value = "  +2.65%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "20.38"
$c.ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +2.16%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "1.47"
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +2.93%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +1.27%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.121"
$c.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -1.37%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "4.76"
$c.ClearFormats()
$ws.Cells.Item(33, 5).Value = "  +2.79%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "4.99"
$c.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  +0.11%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.0655"
$c.ClearFormats()
$ws.Cells.Item(35, 5).Value = "  +4.70%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "6.40"
$c.ClearFormats()
$ws.Cells.Item(36, 5).Value = "  -2.77%  "
$ws.Cells.Item(37, 5).Value = "  +0.23%  "
$ws.Cells.Item(38, 5).Value = "  -3.00%  "
$ws.Cells.Item(39, 5).Value = "  +4.91%  "
$ws.Cells.Item(40, 5).Value = "  -0.02%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "8.67"
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = "  +1.10%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.000225"
$c.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -8.57%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.0961"
$c.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -2.08%  "
$ws.Cells.Item(44, 5).Value = "  +0.87%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "97.22"
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -1.85%  "
$ws.Cells.Item(46, 2).Value = "FTXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "4.38"
$c.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -9.05%  "
$ws.Cells.Item(47, 2).Value = "Maker"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "1.459.53"
$c.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -0.33%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "16.20"
$c.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  -2.23%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "2.74"
$c.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -1.23%  "
$ws.Cells.Item(50, 5).Value = "  -0.73%  "
$ws.Cells.Item(51, 2).Value = "MultiversX"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "51.29"
$c.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -1.93%  "
